# Slide 1, shape "Rectangle 19" currently reads:
#   "GROUP 7 PROJECT PRESENTATION "
# The edit renames the group number: "GROUP 7" -> "GROUP 8", turning the
# single run into two runs ("GROUP 8 " + "PROJECT PRESENTATION ") because
# only the leading "GROUP 7 " characters are retyped.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)

# Locate the title banner shape by name (falls back to the known index
# in case name lookup isn't available in this host).
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Rectangle 19") {
        $shp = $s.Shapes.Item($i)
        break
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(2)
}

$tr  = $shp.TextFrame.TextRange

# Remove the old "GROUP 7 " prefix (first 8 characters) ...
$tr.Characters(1, 8).Text = ""

# ... then type the new prefix back in at the start of the text range,
# leaving "PROJECT PRESENTATION " (and its run/formatting) untouched.
[void]$tr.InsertBefore("GROUP 8 ")
